$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Activation date
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2025"

# 2. Objectives (Portuguese)
Replace-Text "Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II." "Desenvolver e aperfeiçoar o entendimento teórico e prático dos processos bioquímicos fundamentais através da realização de atividades práticas de laboratório."

# 3. Objectives (English) - insert text into previously-empty italic run
$p = $d.Paragraphs.Item(7)
$p.Range.InsertBefore("Developing and enhancing the theoretical and practical understanding of fundamental biochemical processes through the execution of laboratory practical activities.")

# 4. Programa resumido (Portuguese)
Replace-Text "Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas." "Reação de saponificação; Enzimas proteolíticas em produtos comerciais; Extração líquido-líquido de proteínas; Biomateriais sustentáveis; Produção e destilação de etanol; Precipitação de biomoléculas e Reação de Hill."

# 5. Programa resumido (English) - the source text contains stray Greek
#    question-mark characters standing in place of semicolons, so target
#    this paragraph directly instead of relying on Find text matching.
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "Saponification reaction; Proteolytic enzymes in commercial products; Liquid-liquid extraction of proteins; Sustainable biomaterials; Production and distillation of ethanol; Precipitation of biomolecules; and Hill reaction."

# 6. Programa (Portuguese) - contains straight quotes, so assign Range.Text
#    directly rather than Find/Replace to avoid smart-quote autocorrection.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "Aplicação da saponificação em processos industriais, agentes envolvidos na reação e sua utilização em produtos comerciais (cálculo de rendimento, CMC e pH). Avaliação enzimática de proteases para determinação de sua atividade proteolítica em produtos comerciais (sabão em pó, detergentes e cosméticos). Extração líquido-líquido de proteínas e enzimas utilizando solventes orgânicos e polímeros/tensoativos - quantificação dos parâmetros de extração (balanço de massa, recuperação, fator de purificação). Obtenção de biomateriais (bioplástico) de interesse biotecnológico derivado de fontes biológicas - cálculo do rendimento; caracterização do produto final obtido (textura, cor e cheiro) e comparação com os plásticos convencionais. Produção e destilação de etanol - conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono; cálculo da eficiência do processo; ação de um inibidor da glicólise. Precipitação de biomoléculas utilizando diferentes agentes precipitadores (sais, polímeros e solventes orgânicos) - quantificação da recuperação, pH e potencial elétrico. Extração de clorofila e reação de Hill - estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura/luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. *Dentro do programa da disciplina é planejado realizar eventual ""Visita Didática Complementar""."

# 7. Programa (English)
Replace-Text "General proprieties of glycides: main qualitative tests for identification and differentiation of glycides; important reactions and spectrophotometric dosage of reducing monosaccharides.Anaerobic fermentation: general concepts and glucoseFermentation; ethanol and CO2 production; consumption of carbon source; calculation of the process efficiency; action of inhibitor of glycolysis.Chlorophyll extraction andHill reaction: structure of chloroplasts, chlorophyll role in the photosystems I and II; dark phase and light phase; NADP production; dye role as an acceptor ofprotons and electrons. Glycides transportation and enzyme induction: general concepts; enzymes of galactose catabolism; catabolic repression, inactivation and modification; constitutive and induced enzymatic systems in yeast cells." "Application of saponification in industrial processes, agents involved in the reaction, and its use in commercial products (yield calculation, CMC, and pH). Enzymatic evaluation of proteases to determine their proteolytic activity in commercial products (laundry detergent, detergents, and cosmetics). Liquid-liquid extraction of proteins and enzymes using organic solvents and polymers/surfactants - quantification of extraction parameters (mass balance, recovery, purification factor). Production of biomaterials (bioplastic) of biotechnological interest derived from biological sources - yield calculation; characterization of the final product obtained (texture, color, and odor) and comparison with conventional plastics. Production and distillation of ethanol - general concepts and glucose fermentation; ethanol and CO2 production; carbon source consumption; process efficiency calculation; action of a glycolysis inhibitor. Precipitation of biomolecules using different precipitating agents (salts, polymers, and organic solvents) - recovery quantification, pH, and electrical potential. Chlorophyll extraction and Hill reaction - chloroplast structure; role of chlorophyll in photosynthetic systems I and II; dark/light phase; NADP production; ATP production; dye role as proton and electron acceptor. *Complementary didactic visit is planned within the course program."

# 8. Método
Replace-Text "A avaliação será feita por meio de uma prova escrita e notas de relatórios (R)." "A avaliação será realizada através de uma prova escrita (P) e um relatório de atividades práticas (R)."

# 9. Critério
Replace-Text "A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3." "A nota final (NF) será calculada conforme: NF = (P + R)/2. A"

# 10. Norma de recuperação
Replace-Text "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2" "A recuperação será realizada através de uma prova escrita (PR) e a média de recuperação (MR) será calculada conforme: MR = (NF + PR)/2."

# 11. Bibliografia
Replace-Text "CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043." "1. Rocha Filho, J.A., Vitolo, M. Guia para aulas práticas de biotecnologia de enzimas e fermentação. Editora Blucher, 2021. 2. Cisternas, J.R. Fundamentos de bioquímica experimental. São Paulo: Atheneu, 2005. 3. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed Editora, 2022. 4. Voet, D., Voet, J.G., Pratt, C.W. Fundamentos de Bioquímica: a vida em nivel molecular. Artmed Editora, 2014. 5. Vitolo, M., Pessoa Junior, A., Monteiro, G., Carvalho, J.C.M., Stephano, M.A., Sato, S. Biotecnologia farmacêutica: aspectos sobre aplicação industrial. Editora Blucher, 2015."
